# Remove the "Closed": // True, False lines from both objects in the
# "Status" JSON array, drop the now-stray blank paragraph that followed
# the second one, and let Word's own rewrite drop the stale
# lastRenderedPageBreak cached on the following "}, ..." paragraph.

$d = $word.ActiveDocument
$closedLine = "“Closed”: // True, False"

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq $closedLine) {
        # If the paragraph right after this one is just a bare paragraph
        # mark (no real content), it becomes orphaned once "Closed" is
        # gone, so drop it too.
        if ($i -lt $d.Paragraphs.Count) {
            $next = $d.Paragraphs.Item($i + 1)
            if ($next.Range.Text.Length -le 1) {
                $next.Range.Delete()
            }
        }
        $p.Range.Delete()
    }
}

# The paragraph "}, ..." right after the (now last) "Occupied" line in the
# Status array carries a cached <w:lastRenderedPageBreak/> from before the
# edit. Re-assigning its own text forces Word to regenerate the run and
# drop the stale render cache, without altering the visible content.
$closingEllipsis = "}, " + [char]0x2026
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ($closingEllipsis + [char]13)) {
        $p.Range.Text = $closingEllipsis
    }
}
